$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 6021.85
$ws.Range("I28").Value = 839.875
$ws.Range("J28").Value = 26749.75
$ws.Range("K28").Value = 839.875
$ws.Range("L28").Value = 26749.75
$ws.Range("M28").Value = -354.875
$ws.Range("N28").Value = -27719.75
# Row 137
$ws.Range("H137").Value = 3573.4849
$ws.Range("I137").Value = 1749.3158
$ws.Range("J137").Value = 6049.143
$ws.Range("K137").Value = 5247.9474
$ws.Range("L137").Value = 18147.429
$ws.Range("M137").Value = -2697.9474
$ws.Range("N137").Value = -23247.429
# Row 138
$ws.Range("H138").Value = 2908.1428
$ws.Range("J138").Value = 4032.25
$ws.Range("L138").Value = 12096.75
$ws.Range("N138").Value = -22376.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
# Row 26
$ws.Range("H26").Value = 471.14285
$ws.Range("I26").Value = 471.14285
$ws.Range("K26").Value = 471.14285
$ws.Range("M26").Value = -141.14285
# Row 32
$ws.Range("H32").Value = 19341.65
$ws.Range("I32").Value = 16647.182
$ws.Range("J32").Value = 22634.889
$ws.Range("K32").Value = 16647.182
$ws.Range("L32").Value = 22634.889
$ws.Range("M32").Value = -16360.182
$ws.Range("N32").Value = -23208.889
# Row 38
$ws.Range("H38").Value = 8000.5
$ws.Range("I38").Value = 8000.5
$ws.Range("K38").Value = 8000.5
$ws.Range("M38").Value = -7533.5
# Row 39
$ws.Range("H39").Value = 19603.8
$ws.Range("I39").Value = 8750.25
$ws.Range("J39").Value = 63018
$ws.Range("K39").Value = 8750.25
$ws.Range("L39").Value = 63018
$ws.Range("M39").Value = -8230.25
$ws.Range("N39").Value = -64058
# Row 40
$ws.Range("H40").Value = 70000
$ws.Range("J40").Value = 70000
$ws.Range("L40").Value = 70000
$ws.Range("N40").Value = -70352
# Row 42
$ws.Range("H42").Value = 67031
$ws.Range("J42").Value = 67031
$ws.Range("L42").Value = 67031
$ws.Range("N42").Value = -68003
# Row 53
$ws.Range("H53").Value = 25000
$ws.Range("J53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("N53").Value = -26364
# Row 61
$ws.Range("H61").Value = 4452.75
$ws.Range("I61").Value = 4452.75
$ws.Range("K61").Value = 4452.75
$ws.Range("M61").Value = -4240.75
# Row 133
$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -70060
# Row 135
$ws.Range("H135").Value = 94998
$ws.Range("J135").Value = 94998
$ws.Range("L135").Value = 94998
$ws.Range("N135").Value = -105138
# Row 136
$ws.Range("H136").Value = 4452.75
$ws.Range("I136").Value = 4452.75
$ws.Range("K136").Value = 13358.25
$ws.Range("M136").Value = -10808.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 176.14285
$ws.Range("I7").Value = 126.6
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 126.6
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -13.59999999999999
$ws.Range("N7").Value = -526
# Row 22
$ws.Range("H22").Value = 1618.4375
$ws.Range("I22").Value = 433
$ws.Range("J22").Value = 2157.2727
$ws.Range("K22").Value = 433
$ws.Range("L22").Value = 2157.2727
$ws.Range("M22").Value = -83
$ws.Range("N22").Value = -2857.2727
# Row 31
$ws.Range("H31").Value = 5625.069
$ws.Range("I31").Value = 5009.5835
$ws.Range("J31").Value = 8579.4
$ws.Range("K31").Value = 5009.5835
$ws.Range("L31").Value = 8579.4
$ws.Range("M31").Value = -4714.5835
$ws.Range("N31").Value = -9169.4
# Row 34
$ws.Range("H34").Value = 5625.069
$ws.Range("I34").Value = 5009.5835
$ws.Range("J34").Value = 8579.4
$ws.Range("K34").Value = 5009.5835
$ws.Range("L34").Value = 8579.4
$ws.Range("M34").Value = -4807.5835
$ws.Range("N34").Value = -8983.4
# Row 99
$ws.Range("H99").Value = 5549.5557
$ws.Range("I99").Value = 4993.375
$ws.Range("K99").Value = 4993.375
$ws.Range("M99").Value = -3495.375
# Row 105
$ws.Range("H105").Value = 1510.4375
$ws.Range("I105").Value = 1051.6428
$ws.Range("J105").Value = 4722
$ws.Range("K105").Value = 1051.6428
$ws.Range("L105").Value = 4722
$ws.Range("M105").Value = 695.3571999999999
$ws.Range("N105").Value = -8216
# Row 122
$ws.Range("H122").Value = 2477.75
$ws.Range("I122").Value = 2477.75
$ws.Range("K122").Value = 7433.25
$ws.Range("M122").Value = -4983.25
# Row 126
$ws.Range("H126").Value = 5549.5557
$ws.Range("I126").Value = 4993.375
$ws.Range("K126").Value = 14980.125
$ws.Range("M126").Value = -12510.125

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 5206.2
$ws.Range("J114").Value = 5257.75
$ws.Range("L114").Value = 15773.25
$ws.Range("N114").Value = -22281.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 1080.75
$ws.Range("J23").Value = 1080.75
$ws.Range("L23").Value = 1080.75
$ws.Range("N23").Value = -1526.75
# Row 52
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518
# Row 55
$ws.Range("H55").Value = 6937.3335
$ws.Range("I55").Value = 5924.8
$ws.Range("K55").Value = 5924.8
$ws.Range("M55").Value = -5597.8
# Row 92
$ws.Range("H92").Value = 6156.375
$ws.Range("J92").Value = 6535.857
$ws.Range("L92").Value = 6535.857
$ws.Range("N92").Value = -10279.857
# Row 126
$ws.Range("H126").Value = 3436.5
$ws.Range("I126").Value = 3436.5
$ws.Range("K126").Value = 10309.5
$ws.Range("M126").Value = -7839.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1282.3636
$ws.Range("I22").Value = 881.4
$ws.Range("J22").Value = 1616.5
$ws.Range("K22").Value = 881.4
$ws.Range("L22").Value = 1616.5
$ws.Range("M22").Value = -586.4
$ws.Range("N22").Value = -2206.5
# Row 27
$ws.Range("H27").Value = 1282.3636
$ws.Range("I27").Value = 881.4
$ws.Range("J27").Value = 1616.5
$ws.Range("K27").Value = 881.4
$ws.Range("L27").Value = 1616.5
$ws.Range("M27").Value = -774.4
$ws.Range("N27").Value = -1830.5
# Row 40
$ws.Range("H40").Value = 9267.333000000001
$ws.Range("I40").Value = 8204.299999999999
$ws.Range("K40").Value = 8204.299999999999
$ws.Range("M40").Value = -8068.299999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 820
$ws.Range("I113").Value = 900
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = -530
# Row 132
$ws.Range("H132").Value = 5556
$ws.Range("I132").Value = 5148.467
$ws.Range("K132").Value = 15445.401
$ws.Range("M132").Value = -12915.401
